# Progress-Report.xlsx edit:
#  - Paperworks sheet: split the "Class Diagram" paper into two iteration
#    rows. Iteration 1 is now done (100%), Iteration 2 is the new
#    "to be done" item (0%).
#  - Re-point the active tab / selection at the "User Interface" sheet
#    (was left on "Paperworks").

$wb = $excel.ActiveWorkbook

# --- Paperworks sheet: update the Class Diagram rows -----------------
$paperworks = $wb.Worksheets.Item("Paperworks")

# Row 6 was "Class Diagram" / 0% / "To be done during Integration II".
# It becomes "Iteration 1 Class Diagram " and is now fully done, with no
# remaining remark.
$paperworks.Range("B6").Value = "Iteration 1 Class Diagram "
$paperworks.Range("C6").Value = 1
$paperworks.Range("E6").Value = ""

# Row 7 is a brand-new entry for the second iteration's class diagram,
# still outstanding.
$paperworks.Range("B7").Value = "Iteration 2 Class Diagram "
$paperworks.Range("C7").Value = 0

# --- Switch the workbook's active sheet / selection -------------------
# Leave Paperworks' own remembered selection on D17 first (selecting on a
# sheet makes it active, so this has to happen before we move on); then
# activate "User Interface" and select C16 there last, so it ends up as
# the workbook's active sheet.
$paperworks.Range("D17").Select()

$userInterface = $wb.Worksheets.Item("User Interface")
$userInterface.Activate()
$userInterface.Range("C16").Select()
